$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Sample_0199", 1)
    ,@("Sample_0172", 0)
    ,@("Sample_0024", 1)
    ,@("Sample_0133", 0)
    ,@("Sample_0214", 1)
    ,@("Sample_0245", 1)
    ,@("Sample_0123", 0)
    ,@("Sample_0106", 0)
    ,@("Sample_0026", 1)
    ,@("Sample_0229", 0)
    ,@("Sample_0216", 0)
    ,@("Sample_0023", 1)
    ,@("Sample_0221", 1)
    ,@("Sample_0111", 0)
    ,@("Sample_0174", 0)
    ,@("Sample_0173", 1)
    ,@("Sample_0205", 1)
    ,@("Sample_0016", 1)
    ,@("Sample_0128", 1)
    ,@("Sample_0247", 1)
    ,@("Sample_0171", 0)
    ,@("Sample_0151", 1)
    ,@("Sample_0198", 0)
    ,@("Sample_0092", 0)
    ,@("Sample_0143", 0)
    ,@("Sample_0120", 0)
    ,@("Sample_0188", 1)
    ,@("Sample_0233", 1)
    ,@("Sample_0045", 0)
    ,@("Sample_0122", 0)
    ,@("Sample_0140", 1)
    ,@("Sample_0179", 1)
    ,@("Sample_0243", 1)
    ,@("Sample_0035", 1)
    ,@("Sample_0197", 1)
    ,@("Sample_0153", 0)
    ,@("Sample_0107", 0)
    ,@("Sample_0027", 1)
    ,@("Sample_0182", 1)
    ,@("Sample_0134", 1)
    ,@("Sample_0012", 1)
    ,@("Sample_0118", 1)
    ,@("Sample_0196", 1)
    ,@("Sample_0155", 0)
    ,@("Sample_0228", 1)
    ,@("Sample_0183", 1)
    ,@("Sample_0225", 0)
    ,@("Sample_0019", 1)
    ,@("Sample_0062", 1)
    ,@("Sample_0222", 1)
    ,@("Sample_0089", 1)
    ,@("Sample_0119", 0)
    ,@("Sample_0195", 1)
    ,@("Sample_0136", 1)
    ,@("Sample_0178", 0)
    ,@("Sample_0184", 0)
    ,@("Sample_0044", 0)
    ,@("Sample_0192", 1)
    ,@("Sample_0230", 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Ensure the newly added row (60) carries the same cell style as the rows above it.
$ws.Cells.Item(59, 1).Copy() | Out-Null
$ws.Cells.Item(60, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(60, 1).Value = $data[58][0]
$excel.CutCopyMode = 0

